$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.480.00'
$ws.Range("E2").Value = '  +3.47%  '
$ws.Range("D3").Value = '1.606.45'
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("E6").Value = '  +6.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '26.84'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +11.20%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.60'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.249'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.51%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0596'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.23%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0912'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.17%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.838.10'
$ws.Range("E13").Value = '  +3.64%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.602.02'
$ws.Range("E14").Value = '  +3.18%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '29.572.80'
$ws.Range("E15").Value = '  +3.94%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.537'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.33%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.75'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.33%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.60%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.58%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.57'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.86%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0₃0692'
$ws.Range("E21").Value = '  +2.94%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.60%  '
$ws.Range("B24").Value = 'Avalanche'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.44%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.06'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.75%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.25'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.32%  '
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.108'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.79%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.00%  '
$ws.Range("B30").Value = 'BinanceUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0470'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.97%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.39%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.50%  '
$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").Value = '1.437.49'
$ws.Range("E34").Value = '  +3.57%  '
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.36%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.54%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.16%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.51'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.47%  '
$ws.Range("B39").Value = 'HuobiToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.29'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0165'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.89%  '
$ws.Range("B41").Value = 'ImmutableX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.530'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.27%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.93'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.792'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.68%  '
$ws.Range("B45").Value = 'BitcoinSV'
$ws.Range("C45").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '52.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +21.07%  '
$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0466'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.37%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.78'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.56%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.749.16'
$ws.Range("E48").Value = '  +3.66%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.27'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.49%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.78'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.64%  '
$ws.Range("B51").Value = 'WEMIXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.837'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.75%  '
